$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'245.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Formula = "'21.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Formula = "'5.407"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Formula = "'0.05812"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Formula = "'3.377"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Formula = "'6.333"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Formula = "'0.8078"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Formula = "'0.9944"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Formula = "'0.01123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Formula = "'0.1426"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Formula = "'0.07490"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Formula = "'0.03194"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Formula = "'0.03056"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Formula = "'4.194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Formula = "'0.09399"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Formula = "'0.001586"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Formula = "'0.04815"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Formula = "'0.006321"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Formula = "'0.004097"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"
$ws.Range("D21").Formula = "'0.0009977"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Formula = "'3.703"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Formula = "'2.239"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Formula = "'0.0003594"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("D40").Formula = "'0.03888"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Formula = "'0.006531"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Formula = "'0.1073"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Formula = "'0.003001"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Formula = "'0.006419"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Formula = "'0.00005593"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Formula = "'0.3901"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Formula = "'0.1460"
$ws.Range("D48").Style = "Normal"
